$wb = $excel.ActiveWorkbook

# ALC!row113
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value2 = 2474.838
$ws.Cells.Item(113, 9).Value2 = 1972.1428
$ws.Cells.Item(113, 10).Value2 = 2592.1333
$ws.Cells.Item(113, 11).Value2 = 1972.1428
$ws.Cells.Item(113, 12).Value2 = 2592.1333
$ws.Cells.Item(113, 13).Value2 = 1281.8572
$ws.Cells.Item(113, 14).Value2 = -9100.1333

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value2 = 2194.9
$ws.Cells.Item(137, 9).Value2 = 1979.75
$ws.Cells.Item(137, 10).Value2 = 3055.5
$ws.Cells.Item(137, 11).Value2 = 5939.25
$ws.Cells.Item(137, 12).Value2 = 9166.5
$ws.Cells.Item(137, 13).Value2 = -3389.25
$ws.Cells.Item(137, 14).Value2 = -14266.5

# ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value2 = 1269.5625
$ws.Cells.Item(45, 9).Value2 = 1165.2142
$ws.Cells.Item(45, 11).Value2 = 1165.2142
$ws.Cells.Item(45, 13).Value2 = -788.2141999999999

# ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value2 = 1777.9615
$ws.Cells.Item(61, 9).Value2 = 1666.3914
$ws.Cells.Item(61, 10).Value2 = 2633.3333
$ws.Cells.Item(61, 11).Value2 = 1666.3914
$ws.Cells.Item(61, 12).Value2 = 2633.3333
$ws.Cells.Item(61, 13).Value2 = -1454.3914
$ws.Cells.Item(61, 14).Value2 = -3057.3333

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value2 = 7145762.5
$ws.Cells.Item(74, 9).Value2 = 8698187
$ws.Cells.Item(74, 10).Value2 = 4608.4
$ws.Cells.Item(74, 11).Value2 = 8698187
$ws.Cells.Item(74, 12).Value2 = 4608.4
$ws.Cells.Item(74, 13).Value2 = -8697313
$ws.Cells.Item(74, 14).Value2 = -6356.4

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value2 = 7145762.5
$ws.Cells.Item(77, 9).Value2 = 8698187
$ws.Cells.Item(77, 10).Value2 = 4608.4
$ws.Cells.Item(77, 11).Value2 = 43490935
$ws.Cells.Item(77, 12).Value2 = 23042
$ws.Cells.Item(77, 13).Value2 = -43486567
$ws.Cells.Item(77, 14).Value2 = -31778

# ARM!row102
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value2 = 1941.6666
$ws.Cells.Item(102, 9).Value2 = 1776.75
$ws.Cells.Item(102, 11).Value2 = 1776.75
$ws.Cells.Item(102, 13).Value2 = -154.75

# ARM!row110
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value2 = 5641.0435
$ws.Cells.Item(110, 9).Value2 = 6433.8945
$ws.Cells.Item(110, 10).Value2 = 1875
$ws.Cells.Item(110, 11).Value2 = 6433.8945
$ws.Cells.Item(110, 12).Value2 = 1875
$ws.Cells.Item(110, 13).Value2 = -4388.8945
$ws.Cells.Item(110, 14).Value2 = -5965

# ARM!row123
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(123, 8).Value2 = 58266.668
$ws.Cells.Item(123, 10).Value2 = 58266.668
$ws.Cells.Item(123, 12).Value2 = 58266.668
$ws.Cells.Item(123, 14).Value2 = -68066.66800000001

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value2 = 2322.225
$ws.Cells.Item(132, 9).Value2 = 1968.862
$ws.Cells.Item(132, 10).Value2 = 3253.818
$ws.Cells.Item(132, 11).Value2 = 5906.586
$ws.Cells.Item(132, 12).Value2 = 9761.454000000002
$ws.Cells.Item(132, 13).Value2 = -3376.586
$ws.Cells.Item(132, 14).Value2 = -14821.454

# ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value2 = 1777.9615
$ws.Cells.Item(136, 9).Value2 = 1666.3914
$ws.Cells.Item(136, 10).Value2 = 2633.3333
$ws.Cells.Item(136, 11).Value2 = 4999.174199999999
$ws.Cells.Item(136, 12).Value2 = 7899.999899999999
$ws.Cells.Item(136, 13).Value2 = -2449.174199999999
$ws.Cells.Item(136, 14).Value2 = -12999.9999

# BSM!row80
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value2 = 102.25
$ws.Cells.Item(80, 9).Value2 = 66.666664
$ws.Cells.Item(80, 10).Value2 = 123.6
$ws.Cells.Item(80, 11).Value2 = 66.666664
$ws.Cells.Item(80, 12).Value2 = 123.6
$ws.Cells.Item(80, 13).Value2 = 931.333336
$ws.Cells.Item(80, 14).Value2 = -2119.6

# BSM!row83
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(83, 8).Value2 = 102.25
$ws.Cells.Item(83, 9).Value2 = 66.666664
$ws.Cells.Item(83, 10).Value2 = 123.6
$ws.Cells.Item(83, 11).Value2 = 333.33332
$ws.Cells.Item(83, 12).Value2 = 618
$ws.Cells.Item(83, 13).Value2 = 4658.66668
$ws.Cells.Item(83, 14).Value2 = -10602

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value2 = 2950.4736
$ws.Cells.Item(134, 9).Value2 = 2205.4546
$ws.Cells.Item(134, 10).Value2 = 3974.875
$ws.Cells.Item(134, 11).Value2 = 6616.3638
$ws.Cells.Item(134, 12).Value2 = 11924.625
$ws.Cells.Item(134, 13).Value2 = -4081.3638
$ws.Cells.Item(134, 14).Value2 = -16994.625

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value2 = 1792.2
$ws.Cells.Item(58, 9).Value2 = 1312.5555
$ws.Cells.Item(58, 10).Value2 = 3025.5715
$ws.Cells.Item(58, 11).Value2 = 1312.5555
$ws.Cells.Item(58, 12).Value2 = 3025.5715
$ws.Cells.Item(58, 13).Value2 = -1109.5555
$ws.Cells.Item(58, 14).Value2 = -3431.5715

# CRP!row94
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(94, 8).Value2 = 66667708
$ws.Cells.Item(94, 9).Value2 = 142858000
$ws.Cells.Item(94, 10).Value2 = 1201.75
$ws.Cells.Item(94, 11).Value2 = 142858000
$ws.Cells.Item(94, 12).Value2 = 1201.75
$ws.Cells.Item(94, 13).Value2 = -142857549
$ws.Cells.Item(94, 14).Value2 = -2103.75

# CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value2 = 2106177.8
$ws.Cells.Item(99, 9).Value2 = 3251637
$ws.Cells.Item(99, 10).Value2 = 6169
$ws.Cells.Item(99, 11).Value2 = 3251637
$ws.Cells.Item(99, 12).Value2 = 6169
$ws.Cells.Item(99, 13).Value2 = -3250139
$ws.Cells.Item(99, 14).Value2 = -9165

# CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value2 = 2106177.8
$ws.Cells.Item(126, 9).Value2 = 3251637
$ws.Cells.Item(126, 10).Value2 = 6169
$ws.Cells.Item(126, 11).Value2 = 9754911
$ws.Cells.Item(126, 12).Value2 = 18507
$ws.Cells.Item(126, 13).Value2 = -9752441
$ws.Cells.Item(126, 14).Value2 = -23447

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value2 = 11650
$ws.Cells.Item(134, 9).Value2 = 13134
$ws.Cells.Item(134, 11).Value2 = 39402
$ws.Cells.Item(134, 13).Value2 = -36867

# CRP!row135
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(135, 8).Value2 = 6243347.5
$ws.Cells.Item(135, 10).Value2 = 6243347.5
$ws.Cells.Item(135, 12).Value2 = 6243347.5
$ws.Cells.Item(135, 14).Value2 = -6253487.5

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value2 = 1792.2
$ws.Cells.Item(136, 9).Value2 = 1312.5555
$ws.Cells.Item(136, 10).Value2 = 3025.5715
$ws.Cells.Item(136, 11).Value2 = 3937.6665
$ws.Cells.Item(136, 12).Value2 = 9076.7145
$ws.Cells.Item(136, 13).Value2 = -1387.6665
$ws.Cells.Item(136, 14).Value2 = -14176.7145

# CUL!row5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value2 = 1007.88
$ws.Cells.Item(5, 10).Value2 = 1005
$ws.Cells.Item(5, 12).Value2 = 3015
$ws.Cells.Item(5, 14).Value2 = -3239

# CUL!row80
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value2 = 2955.889
$ws.Cells.Item(80, 10).Value2 = 2955.889
$ws.Cells.Item(80, 12).Value2 = 8867.667000000001
$ws.Cells.Item(80, 14).Value2 = -10739.667

# CUL!row83
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(83, 8).Value2 = 2955.889
$ws.Cells.Item(83, 10).Value2 = 2955.889
$ws.Cells.Item(83, 12).Value2 = 26603.001
$ws.Cells.Item(83, 14).Value2 = -35963.001

# CUL!row107
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value2 = 638.4545000000001
$ws.Cells.Item(107, 9).Value2 = 327.16666
$ws.Cells.Item(107, 10).Value2 = 1012
$ws.Cells.Item(107, 11).Value2 = 981.4999799999999
$ws.Cells.Item(107, 12).Value2 = 3036
$ws.Cells.Item(107, 13).Value2 = 938.5000200000001
$ws.Cells.Item(107, 14).Value2 = -6876

# CUL!row113
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value2 = 888.88464
$ws.Cells.Item(113, 9).Value2 = 478.75
$ws.Cells.Item(113, 10).Value2 = 1071.1666
$ws.Cells.Item(113, 11).Value2 = 1436.25
$ws.Cells.Item(113, 12).Value2 = 3213.4998
$ws.Cells.Item(113, 13).Value2 = 733.75
$ws.Cells.Item(113, 14).Value2 = -7553.4998

# CUL!row122
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value2 = 681.0769
$ws.Cells.Item(122, 9).Value2 = 361.875
$ws.Cells.Item(122, 11).Value2 = 3256.875
$ws.Cells.Item(122, 13).Value2 = -806.875

# CUL!row132
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value2 = 1290.8
$ws.Cells.Item(132, 9).Value2 = 834.6667
$ws.Cells.Item(132, 10).Value2 = 1975
$ws.Cells.Item(132, 11).Value2 = 7512.0003
$ws.Cells.Item(132, 12).Value2 = 17775
$ws.Cells.Item(132, 13).Value2 = -4982.0003
$ws.Cells.Item(132, 14).Value2 = -22835

# CUL!row135
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value2 = 1007.88
$ws.Cells.Item(135, 10).Value2 = 1005
$ws.Cells.Item(135, 12).Value2 = 9045
$ws.Cells.Item(135, 14).Value2 = -14115

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value2 = 1531.7778
$ws.Cells.Item(102, 9).Value2 = 1531.7273
$ws.Cells.Item(102, 10).Value2 = 1532
$ws.Cells.Item(102, 11).Value2 = 1531.7273
$ws.Cells.Item(102, 12).Value2 = 1532
$ws.Cells.Item(102, 13).Value2 = 90.27269999999999
$ws.Cells.Item(102, 14).Value2 = -4776

# GSM!row123
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(123, 8).Value2 = 9993.714
$ws.Cells.Item(123, 10).Value2 = 9993.714
$ws.Cells.Item(123, 12).Value2 = 9993.714
$ws.Cells.Item(123, 14).Value2 = -14893.714

# GSM!row126
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value2 = 3312.1904
$ws.Cells.Item(126, 9).Value2 = 3800.0667
$ws.Cells.Item(126, 10).Value2 = 2092.5
$ws.Cells.Item(126, 11).Value2 = 11400.2001
$ws.Cells.Item(126, 12).Value2 = 6277.5
$ws.Cells.Item(126, 13).Value2 = -8930.2001
$ws.Cells.Item(126, 14).Value2 = -11217.5

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value2 = 2782.7273
$ws.Cells.Item(132, 9).Value2 = 1118.6666
$ws.Cells.Item(132, 10).Value2 = 4779.6
$ws.Cells.Item(132, 11).Value2 = 3355.9998
$ws.Cells.Item(132, 12).Value2 = 14338.8
$ws.Cells.Item(132, 13).Value2 = -825.9998000000001
$ws.Cells.Item(132, 14).Value2 = -19398.8

# GSM!row134
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(134, 8).Value2 = 26786.924
$ws.Cells.Item(134, 10).Value2 = 26786.924
$ws.Cells.Item(134, 12).Value2 = 80360.772
$ws.Cells.Item(134, 14).Value2 = -85430.772

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value2 = 1775.5714
$ws.Cells.Item(7, 9).Value2 = 1370.2
$ws.Cells.Item(7, 11).Value2 = 1370.2
$ws.Cells.Item(7, 13).Value2 = -1258.2

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value2 = 1775.5714
$ws.Cells.Item(126, 9).Value2 = 1370.2
$ws.Cells.Item(126, 11).Value2 = 4110.6
$ws.Cells.Item(126, 13).Value2 = -1640.6

# LTW!row135
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(135, 8).Value2 = 50000
$ws.Cells.Item(135, 10).Value2 = 50000
$ws.Cells.Item(135, 12).Value2 = 50000
$ws.Cells.Item(135, 14).Value2 = -60140

# WVR!row15
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value2 = 7000
$ws.Cells.Item(15, 9).Value2 = 0
$ws.Cells.Item(15, 10).Value2 = 7000
$ws.Cells.Item(15, 11).Value2 = 0
$ws.Cells.Item(15, 12).Value2 = 7000
$ws.Cells.Item(15, 14).Value2 = -7576
$ws.Cells.Item(15, 13).ClearContents()

# WVR!row46
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value2 = 0
$ws.Cells.Item(46, 10).Value2 = 0
$ws.Cells.Item(46, 12).Value2 = 0
$ws.Cells.Item(46, 14).ClearContents()

# WVR!row122
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value2 = 100001730
$ws.Cells.Item(122, 9).Value2 = 200001090
$ws.Cells.Item(122, 11).Value2 = 600003270
$ws.Cells.Item(122, 13).Value2 = -600000820

# WVR!row123
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(123, 8).Value2 = 30116.125
$ws.Cells.Item(123, 10).Value2 = 30116.125
$ws.Cells.Item(123, 12).Value2 = 30116.125
$ws.Cells.Item(123, 14).Value2 = -39916.125

# WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value2 = 4601.4546
$ws.Cells.Item(126, 9).Value2 = 5933.25
$ws.Cells.Item(126, 10).Value2 = 1050
$ws.Cells.Item(126, 11).Value2 = 17799.75
$ws.Cells.Item(126, 12).Value2 = 3150
$ws.Cells.Item(126, 13).Value2 = -15329.75
$ws.Cells.Item(126, 14).Value2 = -8090

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value2 = 2934.1765
$ws.Cells.Item(132, 9).Value2 = 2379.1924
$ws.Cells.Item(132, 10).Value2 = 4737.875
$ws.Cells.Item(132, 11).Value2 = 7137.5772
$ws.Cells.Item(132, 12).Value2 = 14213.625
$ws.Cells.Item(132, 13).Value2 = -4607.5772
$ws.Cells.Item(132, 14).Value2 = -19273.625

# WVR!row134
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(134, 8).Value2 = 0
$ws.Cells.Item(134, 10).Value2 = 0
$ws.Cells.Item(134, 12).Value2 = 0
$ws.Cells.Item(134, 14).ClearContents()

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value2 = 1109.05
$ws.Cells.Item(136, 9).Value2 = 833.9729599999999
$ws.Cells.Item(136, 11).Value2 = 2501.91888
$ws.Cells.Item(136, 13).Value2 = 48.08112000000028
